# Fruta / hortaliza, semanal
# A new weekly observation is inserted as row 3, pushing the previously
# existing rows 3-17 down to 4-18 (i.e. a new row is inserted above the
# old row 3, not appended at the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3; this shifts old rows 3..17
# down to 4..18 and carries their formatting (incl. the date style on
# column D) along for the ride.
$ws.Rows(3).Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Range("A3").Value = 9
$ws.Range("B3").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C3").Value = "Metropolitana"
$ws.Range("D3").Value = 45282
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101008
$ws.Range("J3").Value = "Mora"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 400
$ws.Range("N3").Value = 4000
$ws.Range("O3").Value = 4000
$ws.Range("P3").Value = 4000
$ws.Range("Q3").Value = '$/bandeja 2 kilos'
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 2000
$ws.Range("T3").Value = 2
